$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 4740.077
$ws.Range("I12").Value = 2550.125
$ws.Range("K12").Value = 2550.125
$ws.Range("M12").Value = -2380.125
$ws.Range("H33").Value = 179.07143
$ws.Range("I33").Value = 212.75
$ws.Range("J33").Value = 153.8125
$ws.Range("K33").Value = 212.75
$ws.Range("L33").Value = 153.8125
$ws.Range("M33").Value = 16.25
$ws.Range("N33").Value = -611.8125
$ws.Range("H40").Value = 9453.4375
$ws.Range("I40").Value = 7938.25
$ws.Range("K40").Value = 7938.25
$ws.Range("M40").Value = -7763.25
$ws.Range("H58").Value = 2974.25
$ws.Range("J58").Value = 8242.25
$ws.Range("L58").Value = 24726.75
$ws.Range("N58").Value = -25026.75

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 3010.4285
$ws.Range("I61").Value = 1148
$ws.Range("K61").Value = 1148
$ws.Range("M61").Value = -936
$ws.Range("H63").Value = 6799.75
$ws.Range("I63").Value = 6799.75
$ws.Range("K63").Value = 6799.75
$ws.Range("M63").Value = -6113.75
$ws.Range("H66").Value = 6799.75
$ws.Range("I66").Value = 6799.75
$ws.Range("K66").Value = 33998.75
$ws.Range("M66").Value = -30566.75
$ws.Range("H74").Value = 12502.583
$ws.Range("I74").Value = 14656.211
$ws.Range("J74").Value = 4318.8
$ws.Range("K74").Value = 14656.211
$ws.Range("L74").Value = 4318.8
$ws.Range("M74").Value = -13782.211
$ws.Range("N74").Value = -6066.8
$ws.Range("H77").Value = 12502.583
$ws.Range("I77").Value = 14656.211
$ws.Range("J77").Value = 4318.8
$ws.Range("K77").Value = 73281.05499999999
$ws.Range("L77").Value = 21594
$ws.Range("M77").Value = -68913.05499999999
$ws.Range("N77").Value = -30330
$ws.Range("H110").Value = 669085.7
$ws.Range("I110").Value = 859424.5600000001
$ws.Range("J110").Value = 2899.5
$ws.Range("K110").Value = 859424.5600000001
$ws.Range("L110").Value = 2899.5
$ws.Range("M110").Value = -857379.5600000001
$ws.Range("N110").Value = -6989.5
$ws.Range("H122").Value = 3454.6538
$ws.Range("I122").Value = 2300.8125
$ws.Range("J122").Value = 5300.8
$ws.Range("K122").Value = 6902.4375
$ws.Range("L122").Value = 15902.4
$ws.Range("M122").Value = -4452.4375
$ws.Range("N122").Value = -20802.4
$ws.Range("H132").Value = 7686.5483
$ws.Range("I132").Value = 5567.0625
$ws.Range("K132").Value = 16701.1875
$ws.Range("M132").Value = -14171.1875
$ws.Range("H136").Value = 3010.4285
$ws.Range("I136").Value = 1148
$ws.Range("K136").Value = 3444
$ws.Range("M136").Value = -894

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 2481.8572
$ws.Range("I86").Value = 1474.8
$ws.Range("J86").Value = 4999.5
$ws.Range("K86").Value = 1474.8
$ws.Range("L86").Value = 4999.5
$ws.Range("M86").Value = -351.8
$ws.Range("N86").Value = -7245.5
$ws.Range("H89").Value = 2481.8572
$ws.Range("I89").Value = 1474.8
$ws.Range("J89").Value = 4999.5
$ws.Range("K89").Value = 7374
$ws.Range("L89").Value = 24997.5
$ws.Range("M89").Value = -1758
$ws.Range("N89").Value = -36229.5
$ws.Range("H105").Value = 42033.96
$ws.Range("I105").Value = 126489.5
$ws.Range("J105").Value = 2290.1765
$ws.Range("K105").Value = 126489.5
$ws.Range("L105").Value = 2290.1765
$ws.Range("M105").Value = -124742.5
$ws.Range("N105").Value = -5784.1765
$ws.Range("H107").Value = 444
$ws.Range("I107").Value = 320.375
$ws.Range("K107").Value = 320.375
$ws.Range("M107").Value = 1599.625

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H107").Value = 1280.7838
$ws.Range("I107").Value = 983.7917
$ws.Range("J107").Value = 1829.0769
$ws.Range("K107").Value = 983.7917
$ws.Range("L107").Value = 1829.0769
$ws.Range("M107").Value = 936.2083
$ws.Range("N107").Value = -5669.0769

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H86").Value = 5500001.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 5500001.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 16500004.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -16502376.5
$ws.Range("H89").Value = 5500001.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 5500001.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 49500013.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -49511869.5
$ws.Range("H92").Value = 659
$ws.Range("I92").Value = 312.5
$ws.Range("J92").Value = 857
$ws.Range("K92").Value = 937.5
$ws.Range("L92").Value = 2571
$ws.Range("M92").Value = 310.5
$ws.Range("N92").Value = -5067
$ws.Range("H110").Value = 4999.5
$ws.Range("I110").Value = 4999.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 14998.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -10908.5
$ws.Range("N110").ClearContents()
$ws.Range("H128").Value = 279655.66
$ws.Range("I128").Value = 279655.66
$ws.Range("K128").Value = 838966.98
$ws.Range("M128").Value = -833986.98
$ws.Range("H131").Value = 4632.1113
$ws.Range("J131").Value = 7923.778
$ws.Range("L131").Value = 23771.334
$ws.Range("N131").Value = -33851.334
$ws.Range("H134").Value = 7910.7646
$ws.Range("I134").Value = 6391.643
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 19174.929
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -14104.929
$ws.Range("N134").Value = -55140
$ws.Range("H137").Value = 2627.4285
$ws.Range("I137").Value = 2623.6667
$ws.Range("J137").Value = 2650
$ws.Range("K137").Value = 7871.000100000001
$ws.Range("L137").Value = 7950
$ws.Range("M137").Value = -2771.000100000001
$ws.Range("N137").Value = -18150
$ws.Range("H141").Value = 7689.143
$ws.Range("I141").Value = 6764.8
$ws.Range("K141").Value = 20294.4
$ws.Range("M141").Value = -15114.4

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H29").Value = 250000750
$ws.Range("I29").Value = 999.6667
$ws.Range("K29").Value = 999.6667
$ws.Range("M29").Value = -709.6667
$ws.Range("H46").Value = 24247
$ws.Range("I46").Value = 22020.5
$ws.Range("J46").Value = 28700
$ws.Range("K46").Value = 22020.5
$ws.Range("L46").Value = 28700
$ws.Range("M46").Value = -21864.5
$ws.Range("N46").Value = -29012

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 376044.12
$ws.Range("I7").Value = 559734.9
$ws.Range("J7").Value = 8662.556
$ws.Range("K7").Value = 559734.9
$ws.Range("L7").Value = 8662.556
$ws.Range("M7").Value = -559622.9
$ws.Range("N7").Value = -8886.556
$ws.Range("H40").Value = 591856
$ws.Range("I40").Value = 772580
$ws.Range("J40").Value = 4503
$ws.Range("K40").Value = 772580
$ws.Range("L40").Value = 4503
$ws.Range("M40").Value = -772444
$ws.Range("N40").Value = -4775
$ws.Range("H42").Value = 29928
$ws.Range("J42").Value = 29928
$ws.Range("L42").Value = 29928
$ws.Range("N42").Value = -31054
$ws.Range("H49").Value = 29928
$ws.Range("J49").Value = 29928
$ws.Range("L49").Value = 29928
$ws.Range("N49").Value = -30222
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H121").Value = 68624.5
$ws.Range("J121").Value = 68624.5
$ws.Range("L121").Value = 68624.5
$ws.Range("N121").Value = -72118.5
$ws.Range("H122").Value = 640346.2
$ws.Range("I122").Value = 502653.66
$ws.Range("K122").Value = 1507960.98
$ws.Range("M122").Value = -1505510.98
$ws.Range("H126").Value = 376044.12
$ws.Range("I126").Value = 559734.9
$ws.Range("J126").Value = 8662.556
$ws.Range("K126").Value = 1679204.7
$ws.Range("L126").Value = 25987.668
$ws.Range("M126").Value = -1676734.7
$ws.Range("N126").Value = -30927.668
$ws.Range("H132").Value = 4082.6667
$ws.Range("I132").Value = 1783.2858
$ws.Range("J132").Value = 5545.909
$ws.Range("K132").Value = 5349.857400000001
$ws.Range("L132").Value = 16637.727
$ws.Range("M132").Value = -2819.857400000001
$ws.Range("N132").Value = -21697.727
$ws.Range("H136").Value = 2937.4194
$ws.Range("I136").Value = 2217.7693
$ws.Range("K136").Value = 6653.3079
$ws.Range("M136").Value = -4103.3079
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 69535.766
$ws.Range("I62").Value = 115847.336
$ws.Range("J62").Value = 17435.25
$ws.Range("K62").Value = 115847.336
$ws.Range("L62").Value = 17435.25
$ws.Range("M62").Value = -115223.336
$ws.Range("N62").Value = -18683.25
$ws.Range("H65").Value = 69535.766
$ws.Range("I65").Value = 115847.336
$ws.Range("J65").Value = 17435.25
$ws.Range("K65").Value = 579236.6799999999
$ws.Range("L65").Value = 87176.25
$ws.Range("M65").Value = -576116.6799999999
$ws.Range("N65").Value = -93416.25
$ws.Range("H107").Value = 45359.695
$ws.Range("I107").Value = 57605.723
$ws.Range("K107").Value = 172817.169
$ws.Range("M107").Value = -170897.169
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 149720
$ws.Range("J111").Value = 149720
$ws.Range("L111").Value = 149720
$ws.Range("N111").Value = -157900
$ws.Range("H133").Value = 55331.5
$ws.Range("J133").Value = 55331.5
$ws.Range("L133").Value = 55331.5
$ws.Range("N133").Value = -65451.5
